$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.74466443703171
$ws.Range("C2").Value = 9.638021098430626
$ws.Range("D2").Value = 3.635263872756191
$ws.Range("F2").Value = 20.32543466251928
$ws.Range("G2").Value = 22.54005554785133
$ws.Range("H2").Value = 12.28054009667334
$ws.Range("I2").Value = 16.67854648320278
$ws.Range("N2").Value = 15.67933121680821

$ws.Range("B3").Value = 12.0893977438825
$ws.Range("C3").Value = 9.080374938754399
$ws.Range("D3").Value = 3.598498053238155
$ws.Range("F3").Value = 20.23997214437274
$ws.Range("G3").Value = 22.30931844334712
$ws.Range("H3").Value = 12.31745065202186
$ws.Range("I3").Value = 16.78018296761729
$ws.Range("N3").Value = 15.71968328220642

$ws.Range("B4").Value = 11.66940778332275
$ws.Range("C4").Value = 8.718244078861931
$ws.Range("D4").Value = 3.575582800023437
$ws.Range("F4").Value = 20.19585770397381
$ws.Range("G4").Value = 22.17921669237902
$ws.Range("H4").Value = 12.34336436708457
$ws.Range("I4").Value = 16.84772660855484
$ws.Range("N4").Value = 15.74632540445379

$ws.Range("B5").Value = 11.49401176588586
$ws.Range("C5").Value = 8.565754046197336
$ws.Range("D5").Value = 3.566164386655629
$ws.Range("F5").Value = 20.17999620517184
$ws.Range("G5").Value = 22.12917808075074
$ws.Range("H5").Value = 12.35473866293884
$ws.Range("I5").Value = 16.87653751461206
$ws.Range("N5").Value = 15.7576520819536

$ws.Range("B6").Value = 11.46463719642007
$ws.Range("C6").Value = 8.540137616781339
$ws.Range("D6").Value = 3.564595779275534
$ws.Range("F6").Value = 20.17749050013318
$ws.Range("G6").Value = 22.12105096897592
$ws.Range("H6").Value = 12.35667644054389
$ws.Range("I6").Value = 16.88139906679037
$ws.Range("N6").Value = 15.75956126139564

$ws.Range("B7").Value = 11.66705924682604
$ws.Range("C7").Value = 8.716207394796658
$ws.Range("D7").Value = 3.57545609774275
$ws.Range("F7").Value = 20.19563520997981
$ws.Range("G7").Value = 22.17852970806766
$ws.Range("H7").Value = 12.34351447221689
$ws.Range("I7").Value = 16.84810996185164
$ws.Range("N7").Value = 15.74647625696813

$ws.Range("B8").Value = 12.52249779815499
$ws.Range("C8").Value = 9.449863039232413
$ws.Range("D8").Value = 3.622661023714937
$ws.Range("F8").Value = 20.29424126054191
$ws.Range("G8").Value = 22.45813985031164
$ws.Range("H8").Value = 12.29259039090564
$ws.Range("I8").Value = 16.71252077436454
$ws.Range("N8").Value = 15.69285773810003

$ws.Range("B9").Value = 14.05305812566807
$ws.Range("C9").Value = 10.73061101368953
$ws.Range("D9").Value = 3.712289548811381
$ws.Range("F9").Value = 20.55316873383083
$ws.Range("G9").Value = 23.09480525934235
$ws.Range("H9").Value = 12.21865245355457
$ws.Range("I9").Value = 16.48766229576887
$ws.Range("N9").Value = 15.60248884713949

$ws.Range("B10").Value = 15.0806179476739
$ws.Range("C10").Value = 11.57405777934739
$ws.Range("D10").Value = 3.776039179664918
$ws.Range("F10").Value = 20.78208056931259
$ws.Range("G10").Value = 23.61125174877343
$ws.Range("H10").Value = 12.18030621600211
$ws.Range("I10").Value = 16.34783228682168
$ws.Range("N10").Value = 15.54506676037344

$ws.Range("B11").Value = 15.52591551824397
$ws.Range("C11").Value = 11.93646029109731
$ws.Range("D11").Value = 3.804521171068853
$ws.Range("F11").Value = 20.89426868133451
$ws.Range("G11").Value = 23.85554758248814
$ws.Range("H11").Value = 12.16636360030333
$ws.Range("I11").Value = 16.28980815456292
$ws.Range("N11").Value = 15.52088475652916

$ws.Range("B12").Value = 15.69128329291728
$ws.Range("C12").Value = 12.0706266229702
$ws.Range("D12").Value = 3.815226892626812
$ws.Range("F12").Value = 20.93787567580343
$ws.Range("G12").Value = 23.94929459306107
$ws.Range("H12").Value = 12.16158975103367
$ws.Range("I12").Value = 16.26864572919364
$ws.Range("N12").Value = 15.51200603438204

$ws.Range("B13").Value = 15.65581423231649
$ws.Range("C13").Value = 12.04186800519331
$ws.Range("D13").Value = 3.812924855434499
$ws.Range("F13").Value = 20.92843472796377
$ws.Range("G13").Value = 23.92905121541213
$ws.Range("H13").Value = 12.16259534767088
$ws.Range("I13").Value = 16.27316729826335
$ws.Range("N13").Value = 15.51390584739528

$ws.Range("B14").Value = 15.5395860592895
$ws.Range("C14").Value = 11.94755974149666
$ws.Range("D14").Value = 3.805403566170088
$ws.Range("F14").Value = 20.89783394032118
$ws.Range("G14").Value = 23.86323603563421
$ws.Range("H14").Value = 12.16596070036027
$ws.Range("I14").Value = 16.28805083466681
$ws.Range("N14").Value = 15.52014872002508

$ws.Range("B15").Value = 15.46796691577912
$ws.Range("C15").Value = 11.88939357332603
$ws.Range("D15").Value = 3.80078601350647
$ws.Range("F15").Value = 20.87923535224519
$ws.Range("G15").Value = 23.8230802263239
$ws.Range("H15").Value = 12.16808802987504
$ws.Range("I15").Value = 16.29727312924132
$ws.Range("N15").Value = 15.52400891324106

$ws.Range("B16").Value = 15.0510645804989
$ws.Range("C16").Value = 11.54994493319764
$ws.Range("D16").Value = 3.77416693421414
$ws.Range("F16").Value = 20.77490801574631
$ws.Range("G16").Value = 23.59546558791547
$ws.Range("H16").Value = 12.18128806678316
$ws.Range("I16").Value = 16.3517372785112
$ws.Range("N16").Value = 15.54668610334147

$ws.Range("B17").Value = 14.78958056873901
$ws.Range("C17").Value = 11.33624634396774
$ws.Range("D17").Value = 3.757700495595582
$ws.Range("F17").Value = 20.71294602785767
$ws.Range("G17").Value = 23.45815264448116
$ws.Range("H17").Value = 12.19028442727335
$ws.Range("I17").Value = 16.3865846782655
$ws.Range("N17").Value = 15.56109428752469

$ws.Range("B18").Value = 14.6371008494546
$ws.Range("C18").Value = 11.21132919206371
$ws.Range("D18").Value = 3.748180983795316
$ws.Range("F18").Value = 20.67806756260341
$ws.Range("G18").Value = 23.38006175542476
$ws.Range("H18").Value = 12.19578832077891
$ws.Range("I18").Value = 16.40715327352796
$ws.Range("N18").Value = 15.56956408229177

$ws.Range("B19").Value = 14.58511886216519
$ws.Range("C19").Value = 11.16869057967295
$ws.Range("D19").Value = 3.744949682784726
$ws.Range("F19").Value = 20.66638992187504
$ws.Range("G19").Value = 23.35377728681092
$ws.Range("H19").Value = 12.19770834576952
$ws.Range("I19").Value = 16.41420745657921
$ws.Range("N19").Value = 15.57246318246803

$ws.Range("B20").Value = 14.81763193696159
$ws.Range("C20").Value = 11.35920232661442
$ws.Range("D20").Value = 3.759458429502053
$ws.Range("F20").Value = 20.71946352617962
$ws.Range("G20").Value = 23.47267872862628
$ws.Range("H20").Value = 12.18929263672363
$ws.Range("I20").Value = 16.38282069739383
$ws.Range("N20").Value = 15.55954161769285

$ws.Range("B21").Value = 15.57381398544446
$ws.Range("C21").Value = 11.97534365483709
$ws.Range("D21").Value = 3.807614958131197
$ws.Range("F21").Value = 20.90679192314922
$ws.Range("G21").Value = 23.88253484258715
$ws.Range("H21").Value = 12.16495846594472
$ws.Range("I21").Value = 16.28365713931881
$ws.Range("N21").Value = 15.51830748256263

$ws.Range("B22").Value = 16.04901574211813
$ws.Range("C22").Value = 12.36014646095991
$ws.Range("D22").Value = 3.838620601412727
$ws.Range("F22").Value = 21.03575422330725
$ws.Range("G22").Value = 24.15756064917856
$ws.Range("H22").Value = 12.15200451982343
$ws.Range("I22").Value = 16.22357287705083
$ws.Range("N22").Value = 15.49298158454391

$ws.Range("B23").Value = 15.79715112755225
$ws.Range("C23").Value = 12.15640739496201
$ws.Range("D23").Value = 3.822116788099317
$ws.Range("F23").Value = 20.96633902553427
$ws.Range("G23").Value = 24.01015542568259
$ws.Range("H23").Value = 12.15864762686804
$ws.Range("I23").Value = 16.25520635973477
$ws.Range("N23").Value = 15.50635012781197

$ws.Range("B24").Value = 14.80495660164884
$ws.Range("C24").Value = 11.34883033969475
$ws.Range("D24").Value = 3.758663831431929
$ws.Range("F24").Value = 20.7165146441068
$ws.Range("G24").Value = 23.46610882333546
$ws.Range("H24").Value = 12.18973999195352
$ws.Range("I24").Value = 16.38452072911846
$ws.Range("N24").Value = 15.56024299923945

$ws.Range("B25").Value = 13.65558347332276
$ws.Range("C25").Value = 10.40119845359761
$ws.Range("D25").Value = 3.688389604794724
$ws.Range("F25").Value = 20.47622222298324
$ws.Range("G25").Value = 22.91362101579332
$ws.Range("H25").Value = 12.23586079531475
$ws.Range("I25").Value = 16.54406289934862
$ws.Range("N25").Value = 15.62535789394857
